# Append 4 new daily rows (17->21 apr / serial 44304-44307) to the
# "Medolla" report sheet, mirroring the existing data layout:
#   A: date serial (formatted like the rows above it)
#   B: nuovi positivi
#   C: somma mobile 7gg.
#   D: somma mobile 7gg. per 100mila abitanti

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 229
$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

$newRows = @(
    @{ Date = 44304; B = 5; C = 12; D = 192.2768787053357 },
    @{ Date = 44305; B = 1; C = 12; D = 192.2768787053357 },
    @{ Date = 44306; B = 0; C = 12; D = 192.2768787053357 },
    @{ Date = 44307; B = 0; C = 12; D = 192.2768787053357 }
)

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $row = $lastRow + 1 + $i
    $data = $newRows[$i]

    $ws.Cells.Item($row, 1).Value = $data.Date
    $ws.Cells.Item($row, 2).Value = $data.B
    $ws.Cells.Item($row, 3).Value = $data.C
    $ws.Cells.Item($row, 4).Value = $data.D

    # Match column A's style (date number format, bold centered, border)
    # used throughout the rest of the sheet.
    $ws.Cells.Item($lastRow, 1).Copy()
    $ws.Cells.Item($row, 1).PasteSpecial($xlPasteFormats)
}

$excel.CutCopyMode = $false
